$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '75.931.37'
$ws.Range("E2").Value = '  +1.90%  '
$ws.Range("D3").Value = '2.914.65'
$ws.Range("E3").Value = '  +4.44%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'" + '199.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.80%  '
$ws.Range("D6").Value = "'" + '598.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = "'" + '0.550'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = "'" + '0.199'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.26%  '
$ws.Range("D10").Value = '2.913.09'
$ws.Range("E10").Value = '  +4.45%  '
$ws.Range("D11").Value = "'" + '0.439'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +19.73%  '
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").Value = "'" + '4.89'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.34%  '
$ws.Range("D14").Value = '3.450.33'
$ws.Range("E14").Value = '  +3.71%  '
$ws.Range("D15").Value = '75.805.85'
$ws.Range("E15").Value = '  +1.71%  '
$ws.Range("D16").Value = "'" + '0.0000191'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.06%  '
$ws.Range("D17").Value = "'" + '27.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.54%  '
$ws.Range("D18").Value = '2.913.88'
$ws.Range("E18").Value = '  +3.85%  '
$ws.Range("D19").Value = "'" + '8.95'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("D20").Value = "'" + '12.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.73%  '
$ws.Range("D21").Value = "'" + '378.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.97%  '
$ws.Range("E22").Value = '  +3.91%  '
$ws.Range("D23").Value = "'" + '4.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.09%  '
$ws.Range("D24").Value = "'" + '71.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.54%  '
$ws.Range("D25").Value = "'" + '0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("D26").Value = '3.065.16'
$ws.Range("E26").Value = '  +3.66%  '
$ws.Range("E27").Value = '  +1.94%  '
$ws.Range("D28").Value = "'" + '9.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.23%  '
$ws.Range("E29").Value = '  +7.30%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("E31").Value = '  +3.31%  '
$ws.Range("D32").Value = "'" + '504.79'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.75%  '
$ws.Range("D33").Value = "'" + '7.73'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("E34").Value = '  +2.46%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = "'" + '164.69'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.94%  '
$ws.Range("D37").Value = "'" + '20.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.77%  '
$ws.Range("B38").Value = 'WhiteBITCoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D38").Value = "'" + '19.67'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.01%  '
$ws.Range("B39").Value = 'Cronos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D39").Value = "'" + '0.106'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +23.52%  '
$ws.Range("E40").Value = '  -3.60%  '
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").Value = "'" + '180.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").Value = "'" + '0.345'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.73%  '
$ws.Range("D44").Value = "'" + '5.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("D46").Value = "'" + '40.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.58%  '
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("E49").Value = '  +2.00%  '
$ws.Range("E50").Value = '  +8.51%  '
$ws.Range("E51").Value = '  +0.59%  '
